$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
    3  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    4  = @(0.3048080303191223, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 2.626907116734944)
    5  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 1, 13.71653804550039)
    6  = @(0.04763786555579896, 0.04240448674262143, 3.900430680208489, 0.496779210170732, 0, 4.487252242677641)
    7  = @(0.3048080303191223, 0.04240448674262143, 0.8054896365839992, 8.660232485948974, 0, 9.812934639594717)
    8  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    9  = @(0.01514828764759746, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 1, 1.35982162114495)
    10 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
